$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (year 2025) metrics
$ws.Range("C6").Value = 363
$ws.Range("D6").Value = 287
$ws.Range("E6").Value = 76
$ws.Range("F6").Value = 63.91982182628062
$ws.Range("G6").Value = 20.9366391184573
$ws.Range("H6").Value = 79.06336088154269
